$wb = $excel.ActiveWorkbook

# --- "final" sheet: header Q1 changes from "ingroup" to the new "rule_J" label ---
$final = $wb.Worksheets.Item("final")
$final.Range("Q1").Value = "rule_J"

# --- "final" sheet: current selection moved to U12 ---
$final.Activate()
$null = $final.Range("U12").Select()

# --- "prep" sheet: the two adjoining "equal to 0" highlight rules (R2:U81 and
#     V2:W81) get consolidated into a single rule covering the combined
#     range R2:W81 ---
$prep = $wb.Worksheets.Item("prep")
$fcKeep = $prep.Range("V2:W81").FormatConditions.Item(1)
$fcKeep.ModifyAppliesToRange($prep.Range("R2:W81"))
$prep.Range("R2:U81").FormatConditions.Delete()
$fcKeep.Priority = 5
